$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 13 (shifts old rows 13-23 down to 15-25)
$ws.Rows("13:14").Insert()

# The freshly-inserted rows only carry column A's bold style by default;
# clear that stray formatting from column A (those rows have no A value)
# and copy the wrap-text body styles from row 10 (B/C) onto the new B/C cells.
$ws.Range("A13:A14").Clear()
$ws.Range("B10").Copy()
$ws.Range("B13:B14").PasteSpecial(-4122)
$ws.Range("C10").Copy()
$ws.Range("C13:C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fix row 10 (Objetivos) - was holding a misplaced docente name, now holds real objectives text
$ws.Range("B10").Value = 'Essa disciplina faz parte da formação do engenheiro de materiais, contribuindo para gerar competências gerais e específicas.Interrelacionar essa disciplina com outras da grade do curso, como: Pirometalurgia, Processamento de Cerâmicas, entre outras. Desenvolver nos alunos a prática da redação científica, da busca bibliográfica, de informações técnicas e especificar equipamentos.Incentivar os alunos ao aproveitamento racional dos recursos naturais transformando-os em produtos com valor agregado e/ou metais, via processo em fase aquosa e temperaturas amenas, levando em consideração os aspectos ambientais, econômicos e sociais.Incentivar trabalhos em grupo, com apresentação de resultados.'
$ws.Range("C10").Value = 'Essa disciplina faz parte da formação do engenheiro de materiais, contribuindo para gerar competências gerais e específicas.Interrelacionar essa disciplina com outras da grade do curso, como: Pirometalurgia, Processamento de Cerâmicas, entre outras. Desenvolver nos alunos a prática da redação científica, da busca bibliográfica, de informações técnicas e especificar equipamentos.Incentivar os alunos ao aproveitamento racional dos recursos naturais transformando-os em produtos com valor agregado e/ou metais, via processo em fase aquosa e temperaturas amenas, levando em consideração os aspectos ambientais, econômicos e sociais.Incentivar trabalhos em grupo, com apresentação de resultados.'

# New row 13: first docente responsavel (Carlos Angelo Nunes)
$ws.Range("B13").Value = '3577649 - Carlos Angelo Nunes'
$ws.Range("C13").Value = '3577649 - Carlos Angelo Nunes'

# New row 14: second docente responsavel (Sebastiao Ribeiro)
$ws.Range("B14").Value = '1922320 - Sebastiao Ribeiro'
$ws.Range("C14").Value = '1922320 - Sebastiao Ribeiro'

# Row 15 (Programa resumido) - was holding a wrong duplicate, now holds real summary text
$ws.Range("B15").Value = '1 – Introdução geral, 2 – Análise granulométrica, 3 – Fragmentação (britagem e moagem), 4 - Classificação e Peneiramento Industrial, 5 – Concentração de minérios, 6 - Circuitos de tratamentos de minérios e balanço de massa, 7 – Matérias-Primas mais usadas em hidrometalurgia, 8 – Lixiviação, 9 - Técnicas para purificação e concentração de licores, 10 – Precipitação, 11 - Circuitos de plantas hidrometalúrgicas e balanço de massa. 12 – Hidrometalurgia e reciclagem, 13 – Testes experimentais, 14 – Visitas técnicas em empresas do ramo, 15 – Em todos os itens anteriores serão abordados aspectos sociais, ambientais, legais e econômicos.'
$ws.Range("C15").Value = '1 – Introdução geral, 2 – Análise granulométrica, 3 – Fragmentação (britagem e moagem), 4 - Classificação e Peneiramento Industrial, 5 – Concentração de minérios, 6 - Circuitos de tratamentos de minérios e balanço de massa, 7 – Matérias-Primas mais usadas em hidrometalurgia, 8 – Lixiviação, 9 - Técnicas para purificação e concentração de licores, 10 – Precipitação, 11 - Circuitos de plantas hidrometalúrgicas e balanço de massa. 12 – Hidrometalurgia e reciclagem, 13 – Testes experimentais, 14 – Visitas técnicas em empresas do ramo, 15 – Em todos os itens anteriores serão abordados aspectos sociais, ambientais, legais e econômicos.'

# Row 17 (Programa) - was holding a wrong duplicate, now holds real full programme text
$ws.Range("B17").Value = '1- Introdução geral sobre Tratamento de Minérios e Hidrometalurgia – Conceitos e terminologias aplicados ao tratamento de minérios e hidrometalurgia, 2 – Análise granulométrica – peneiramento e instrumental, tamanho e distribuição de tamanhos de partículas, 3 – Fragmentação – desmonte, britagem e moagem (teorias, métodos e equipamentos), 4 - Classificação e Peneiramento Industrial – fundamentos da classificação, tipos de classificadores, peneiramento, 5 - Concentração de minérios pelas técnicas de:  sedimentação por queda livre e retardada, gravítica, separação magnética e eletrostática, flotação (físico-química de superfícies; flotação em células e em colunas, um estudo de caso, 6 - Circuitos de tratamento de minérios e balanço de massa, 7 - Matérias-primas mais usadas em hidrometalurgia – minérios e concentrados dos principais metais, materiais secundários, produtos da mineração urbana, 8 – Lixiviação - princípios, técnicas e equipamentos, 9 - Técnicas para purificação e concentração de licores, 10 – Precipitação - princípios, métodos e obtenção de compostos e/ou metais, 11 - Circuitos de plantas hidrometalúrgicas industriais, com ênfase nos metais cobre, níquel, zinco, terras raras, nióbio, tântalo, zircônio e háfnio, envolvendo balanço de massa simples, 12 - Hidrometalurgia como solução para reciclagem e seus efeitos econômicos e sociais, 13 - Testes experimentais; 14 – Visitas técnicas em empresas do ramo. 15 -  Em todos os itens anteriores serão abordados aspectos sociais, ambientais, legais e econômicos para ampliar competências nos alunos.'
$ws.Range("C17").Value = '1- Introdução geral sobre Tratamento de Minérios e Hidrometalurgia – Conceitos e terminologias aplicados ao tratamento de minérios e hidrometalurgia, 2 – Análise granulométrica – peneiramento e instrumental, tamanho e distribuição de tamanhos de partículas, 3 – Fragmentação – desmonte, britagem e moagem (teorias, métodos e equipamentos), 4 - Classificação e Peneiramento Industrial – fundamentos da classificação, tipos de classificadores, peneiramento, 5 - Concentração de minérios pelas técnicas de:  sedimentação por queda livre e retardada, gravítica, separação magnética e eletrostática, flotação (físico-química de superfícies; flotação em células e em colunas, um estudo de caso, 6 - Circuitos de tratamento de minérios e balanço de massa, 7 - Matérias-primas mais usadas em hidrometalurgia – minérios e concentrados dos principais metais, materiais secundários, produtos da mineração urbana, 8 – Lixiviação - princípios, técnicas e equipamentos, 9 - Técnicas para purificação e concentração de licores, 10 – Precipitação - princípios, métodos e obtenção de compostos e/ou metais, 11 - Circuitos de plantas hidrometalúrgicas industriais, com ênfase nos metais cobre, níquel, zinco, terras raras, nióbio, tântalo, zircônio e háfnio, envolvendo balanço de massa simples, 12 - Hidrometalurgia como solução para reciclagem e seus efeitos econômicos e sociais, 13 - Testes experimentais; 14 – Visitas técnicas em empresas do ramo. 15 -  Em todos os itens anteriores serão abordados aspectos sociais, ambientais, legais e econômicos para ampliar competências nos alunos.'

# Row 20 (Metodo) - was holding a misplaced docente name, now holds real method text
$ws.Range("B20").Value = 'Serão realizadas duas provas escritas (P1 e P2) com peso 1. No mínimo, um relatório a partir de trabalhos em grupo, com peso 1 (NR) e avaliação individual realizada durante todo o curso (AI), com peso 1.'
$ws.Range("C20").Value = 'Serão realizadas duas provas escritas (P1 e P2) com peso 1. No mínimo, um relatório a partir de trabalhos em grupo, com peso 1 (NR) e avaliação individual realizada durante todo o curso (AI), com peso 1.'

# Row 21 (Criterio) - was holding the method text (off by one), now holds real criteria text
$ws.Range("B21").Value = 'A nota final (NF) será calculada pela equação: NF = 0,5[(P1 + P2)/2] + 0,3NR + 0,2AI.'
$ws.Range("C21").Value = 'A nota final (NF) será calculada pela equação: NF = 0,5[(P1 + P2)/2] + 0,3NR + 0,2AI.'

# Row 22 (Norma de recuperacao) - was holding the criteria text (off by one), now holds real recovery norm text
$ws.Range("B22").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, inclusive com cobrança das competências desenvolvidas, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'
$ws.Range("C22").Value = 'Para a recuperação será realizada uma prova (PR) abrangendo toda a matéria lecionada no semestre, inclusive com cobrança das competências desenvolvidas, valendo de 0 (zero) a 10 (dez). Média final = (NF + PR)/2.'

# Row 23 (Bibliografia) - was holding the recovery norm text (off by one), now holds real bibliography text
$ws.Range("B23").Value = '1 - Tratamento de Minérios, 5ª Ed., CETEM_MCT, Rio de Janeiro, 2010, Adão Benvindo da Luz, João Alves Sampaio e Salvador L. M. de Almeida, 2 - Dispersão e empacotamento de partículas, Fazendo Arte Editorial, Ivone R. de Oliveira e co-autores, 2.000, 3 - Teoria e prática de tratamento e recuperação de minérios por sistemas gravíticos, Cristoni, S. Signus - São Paulo, 1986, 4 - Teoria e Prática do Tratamento de Minérios – 3ª Edição, Vol. 1-SIGNUS, 5 - Introdução ao Tratamento de Minérios, George Eduardo Sales Valadao, Editora: UFMG, 6 - Principles of extractive metallurgy, Habashi, F., Laval University, vol. 1, Quebec City, Canadá, 1970, 7 - Principles of extractive metallurgy, Habashi, F., Laval University, vol. 2, Quebec City, Canadá, 1970, 8 - Solvent extraction principles and applications to process metal, Ltcey, G.M. & Ashbrook, A.W., Part I – Elsevier, 1984, 9 - Ion exchange resins, Robert Kunin, Robert E. Krieger Publishing Company Malabar, Florida, 1958, 10 – artigos  especializados em Tratamento de Minérios e Hidrometalurgia.'
$ws.Range("C23").Value = '1 - Tratamento de Minérios, 5ª Ed., CETEM_MCT, Rio de Janeiro, 2010, Adão Benvindo da Luz, João Alves Sampaio e Salvador L. M. de Almeida, 2 - Dispersão e empacotamento de partículas, Fazendo Arte Editorial, Ivone R. de Oliveira e co-autores, 2.000, 3 - Teoria e prática de tratamento e recuperação de minérios por sistemas gravíticos, Cristoni, S. Signus - São Paulo, 1986, 4 - Teoria e Prática do Tratamento de Minérios – 3ª Edição, Vol. 1-SIGNUS, 5 - Introdução ao Tratamento de Minérios, George Eduardo Sales Valadao, Editora: UFMG, 6 - Principles of extractive metallurgy, Habashi, F., Laval University, vol. 1, Quebec City, Canadá, 1970, 7 - Principles of extractive metallurgy, Habashi, F., Laval University, vol. 2, Quebec City, Canadá, 1970, 8 - Solvent extraction principles and applications to process metal, Ltcey, G.M. & Ashbrook, A.W., Part I – Elsevier, 1984, 9 - Ion exchange resins, Robert Kunin, Robert E. Krieger Publishing Company Malabar, Florida, 1958, 10 – artigos  especializados em Tratamento de Minérios e Hidrometalurgia.'

# Column B now gets its own explicit width definition, splitting it off from column A's range
$ws.Columns("B").ColumnWidth = 60.7109375
